$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update timestamp text
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 18:24"

# Update country rankings / stats per diff
$ws.Range("B4").Value = 7688210
$ws.Range("C4").Value = 8566
$ws.Range("D4").Value = 4908006
$ws.Range("E4").Value = 2564979
$ws.Range("G4").Value = 193
$ws.Range("H4").Value = 215225
$ws.Range("B5").Value = 6724380
$ws.Range("C5").Value = 42307
$ws.Range("D5").Value = 5703607
$ws.Range("E5").Value = 916741
$ws.Range("G5").Value = 432
$ws.Range("H5").Value = 104032
$ws.Range("B6").Value = 4946913
$ws.Range("C6").Value = 6414
$ws.Range("E6").Value = 504703
$ws.Range("G6").Value = 135
$ws.Range("H6").Value = 146908
$ws.Range("B15").Value = 530113
$ws.Range("C15").Value = 14542
$ws.Range("G15").Value = 76
$ws.Range("H15").Value = 42445
$ws.Range("B31").Value = 142056
$ws.Range("C31").Value = 717
$ws.Range("E31").Value = 9843
$ws.Range("G31").Value = 21
$ws.Range("H31").Value = 11702
$ws.Range("A48").Value = "Chequia"
$ws.Range("B48").Value = 87176
$ws.Range("C48").Value = 1610
$ws.Range("D48").Value = 48714
$ws.Range("E48").Value = 37673
$ws.Range("G48").Value = 31
$ws.Range("H48").Value = 789
$ws.Range("A49").Value = "Japon"
$ws.Range("B49").Value = 86047
$ws.Range("C49").Value = 308
$ws.Range("D49").Value = 79071
$ws.Range("E49").Value = 5374
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 1602
$ws.Range("D60").Value = 57612
$ws.Range("E60").Value = 191
$ws.Range("B86").Value = 20541
$ws.Range("C86").Value = 399
$ws.Range("E86").Value = 10132
$ws.Range("G86").Value = 3
$ws.Range("H86").Value = 420
$ws.Range("A90").Value = "Jordania"
$ws.Range("B90").Value = 19001
$ws.Range("C90").Value = 1537
$ws.Range("D90").Value = 5386
$ws.Range("E90").Value = 13493
$ws.Range("G90").Value = 12
$ws.Range("H90").Value = 122
$ws.Range("A91").Value = "Croacia"
$ws.Range("B91").Value = 18084
$ws.Range("C91").Value = 287
$ws.Range("D91").Value = 16192
$ws.Range("E91").Value = 1588
$ws.Range("G91").Value = 4
$ws.Range("H91").Value = 304
$ws.Range("B96").Value = 14568
$ws.Range("C96").Value = 158
$ws.Range("D96").Value = 8965
$ws.Range("E96").Value = 5200
$ws.Range("G96").Value = 3
$ws.Range("H96").Value = 403
$ws.Range("B100").Value = 12584
$ws.Range("C100").Value = 225
$ws.Range("D100").Value = 8557
$ws.Range("E100").Value = 3839
$ws.Range("G100").Value = 6
$ws.Range("H100").Value = 188
$ws.Range("B106").Value = 10070
$ws.Range("C106").Value = 13
$ws.Range("D106").Value = 9741
$ws.Range("E106").Value = 260
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 69
$ws.Range("B108").Value = 9398
$ws.Range("C108").Value = 102
$ws.Range("D108").Value = 6358
$ws.Range("E108").Value = 2973
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 67
$ws.Range("B111").Value = 8979
$ws.Range("C111").Value = 54
$ws.Range("E111").Value = 1058
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 128
$ws.Range("B149").Value = 2734
$ws.Range("C149").Value = 8
$ws.Range("E149").Value = 1394
$ws.Range("A160").Value = "Republica de Chipre"
$ws.Range("B160").Value = 1876
$ws.Range("C160").Value = 29
$ws.Range("D160").Value = 1369
$ws.Range("E160").Value = 484
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = 23
$ws.Range("A161").Value = "Togo"
$ws.Range("B161").Value = 1864
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 1403
$ws.Range("E161").Value = 413
$ws.Range("H161").Value = 48
$ws.Range("A162").Value = "Nueva Zelanda"
$ws.Range("B162").Value = 1858
$ws.Range("C162").Value = 3
$ws.Range("D162").Value = 1790
$ws.Range("E162").Value = 43
$ws.Range("H162").Value = 25
$ws.Range("D165").Value = 1240
$ws.Range("E165").Value = 32
$ws.Range("B185").Value = 345
$ws.Range("C185").Value = 1
$ws.Range("E185").Value = 5
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Nueva Caledonia"